$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Question/Answer text that changed in the source Q&A table.
# (credit/debit card -> ATM card; answer for "destination" question split
# into two separate answers; "logs" -> "log" typo fix)
$ws.Range("D13").Value = "Destination, transaction log, error log, customer feedbacks, information related to banking systems and digital wallets"
$ws.Range("D4").Value = "ATM card, QR code payment linked with banking system or digital wallet."
$ws.Range("C8").Value = "What happens when a user select ATM card payment option?"
$ws.Range("D5").Value = "The TVM displays a menu of available destinations together with the price to the user."
$ws.Range("D6").Value = "The user chooses the quantity of tickets they want to purchase."

# Match the author's final cursor position/selection in the saved file.
$ws.Range("G6").Select() | Out-Null
